# Daily attendance processing - 2026-01-05 01:45:53
# Swap the order of the "Recorded By" collaborators in column G:
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
# Only cells that contain exactly that combined value are affected;
# cells that list just "System" or just "dnasr281@gmail.com" are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("G:G").Replace(
    "System, dnasr281@gmail.com",
    "dnasr281@gmail.com, System",
    [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole
)
